$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 8818.817999999999
$ws.Range("J51").Value = 9835.286
$ws.Range("L51").Value = 9835.286
$ws.Range("N51").Value = -10803.286
# Row 98
$ws.Range("H98").Value = 838
$ws.Range("I98").Value = 606.95
$ws.Range("K98").Value = 606.95
$ws.Range("M98").Value = 891.05
# Row 122
$ws.Range("H122").Value = 838
$ws.Range("I122").Value = 606.95
$ws.Range("K122").Value = 1820.85
$ws.Range("M122").Value = 629.1499999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 39
$ws.Range("H39").Value = 34000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 34000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = $null
$ws.Range("M39").Value = $null
$ws.Range("N39").Value = -35040
# Row 43
$ws.Range("H43").Value = 18817.25
$ws.Range("I43").Value = 18423.334
$ws.Range("K43").Value = 18423.334
$ws.Range("M43").Value = -18110.334
# Row 102
$ws.Range("H102").Value = 2001
$ws.Range("I102").Value = 1853
$ws.Range("K102").Value = 1853
$ws.Range("M102").Value = -231
# Row 109
$ws.Range("H109").Value = 199483.33
$ws.Range("J109").Value = 199483.33
$ws.Range("L109").Value = 199483.33
$ws.Range("N109").Value = -202257.33
# Row 110
$ws.Range("H110").Value = 2158.875
$ws.Range("I110").Value = 2247.2593
$ws.Range("K110").Value = 2247.2593
$ws.Range("M110").Value = -202.2593000000002
# Row 122
$ws.Range("H122").Value = 3084.9333
$ws.Range("I122").Value = 1881.96
$ws.Range("K122").Value = 5645.88
$ws.Range("M122").Value = -3195.88

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 29
$ws.Range("H29").Value = 3255
$ws.Range("I29").Value = 3255
$ws.Range("K29").Value = 3255
$ws.Range("M29").Value = -2966
# Row 54
$ws.Range("H54").Value = 35999.25
$ws.Range("J54").Value = 35999.25
$ws.Range("L54").Value = 35999.25
$ws.Range("N54").Value = -36967.25
# Row 94
$ws.Range("H94").Value = 709.6857
$ws.Range("I94").Value = 593.76666
$ws.Range("K94").Value = 593.76666
$ws.Range("M94").Value = -142.76666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 50
$ws.Range("H50").Value = 70498
$ws.Range("I50").Value = 69997
$ws.Range("J50").Value = 70598.2
$ws.Range("K50").Value = 69997
$ws.Range("L50").Value = 70598.2
$ws.Range("M50").Value = -69372
$ws.Range("N50").Value = -71848.2
# Row 51
$ws.Range("H51").Value = 58749.75
$ws.Range("J51").Value = 58749.75
$ws.Range("L51").Value = 58749.75
$ws.Range("N51").Value = -60221.75
# Row 61
$ws.Range("H61").Value = 58749.75
$ws.Range("J61").Value = 58749.75
$ws.Range("L61").Value = 58749.75
$ws.Range("N61").Value = -59445.75
# Row 64
$ws.Range("H64").Value = 132601.44
$ws.Range("J64").Value = 132601.44
$ws.Range("L64").Value = 132601.44
$ws.Range("N64").Value = -133097.44
# Row 67
$ws.Range("H67").Value = 132601.44
$ws.Range("J67").Value = 132601.44
$ws.Range("L67").Value = 132601.44
$ws.Range("N67").Value = -134317.44
# Row 94
$ws.Range("H94").Value = 3035.8333
$ws.Range("I94").Value = 1970.8889
$ws.Range("J94").Value = 4100.778
$ws.Range("K94").Value = 1970.8889
$ws.Range("L94").Value = 4100.778
$ws.Range("M94").Value = -1519.8889
$ws.Range("N94").Value = -5002.778
# Row 134
$ws.Range("H134").Value = 3850.8462
$ws.Range("I134").Value = 2795.9524
$ws.Range("K134").Value = 8387.8572
$ws.Range("M134").Value = -5852.8572

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 3714.4546
$ws.Range("I75").Value = 2086
$ws.Range("K75").Value = 6258
$ws.Range("M75").Value = -5260
# Row 78
$ws.Range("H78").Value = 3714.4546
$ws.Range("I78").Value = 2086
$ws.Range("K78").Value = 18774
$ws.Range("M78").Value = -13782
# Row 113
$ws.Range("H113").Value = 606.1429000000001
$ws.Range("J113").Value = 606.1429000000001
$ws.Range("L113").Value = 1818.4287
$ws.Range("N113").Value = -6158.4287
# Row 131
$ws.Range("H131").Value = 792935.0600000001
$ws.Range("J131").Value = 2025137.9
$ws.Range("L131").Value = 6075413.699999999
$ws.Range("N131").Value = -6085493.699999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 9991.866
$ws.Range("I80").Value = 8456.526
$ws.Range("K80").Value = 8456.526
$ws.Range("M80").Value = -7458.526
# Row 83
$ws.Range("H83").Value = 9991.866
$ws.Range("I83").Value = 8456.526
$ws.Range("K83").Value = 42282.63
$ws.Range("M83").Value = -37290.63
# Row 109
$ws.Range("H109").Value = 199500
$ws.Range("J109").Value = 199500
$ws.Range("L109").Value = 199500
$ws.Range("N109").Value = -201580

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 44279.324
$ws.Range("I7").Value = 47931.89
$ws.Range("K7").Value = 47931.89
$ws.Range("M7").Value = -47819.89
# Row 16
$ws.Range("H16").Value = 2720.9524
$ws.Range("I16").Value = 2468.2144
$ws.Range("J16").Value = 3226.4285
$ws.Range("K16").Value = 2468.2144
$ws.Range("L16").Value = 3226.4285
$ws.Range("M16").Value = -2298.2144
$ws.Range("N16").Value = -3566.4285
# Row 46
$ws.Range("H46").Value = 3138.4
$ws.Range("I46").Value = 846.5
$ws.Range("J46").Value = 6576.25
$ws.Range("K46").Value = 846.5
$ws.Range("L46").Value = 6576.25
$ws.Range("M46").Value = -658.5
$ws.Range("N46").Value = -6952.25
# Row 100
$ws.Range("H100").Value = 2462.2
$ws.Range("I100").Value = 2363
$ws.Range("J100").Value = 2611
$ws.Range("K100").Value = 2363
$ws.Range("L100").Value = 2611
$ws.Range("M100").Value = -1822
$ws.Range("N100").Value = -3693
# Row 122
$ws.Range("H122").Value = 4788.923
$ws.Range("I122").Value = 4072.3872
$ws.Range("K122").Value = 12217.1616
$ws.Range("M122").Value = -9767.161599999999
# Row 126
$ws.Range("H126").Value = 44279.324
$ws.Range("I126").Value = 47931.89
$ws.Range("K126").Value = 143795.67
$ws.Range("M126").Value = -141325.67

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 52
$ws.Range("H52").Value = 4299.5
$ws.Range("I52").Value = 4299.5
$ws.Range("K52").Value = 4299.5
$ws.Range("M52").Value = -4073.5
# Row 62
$ws.Range("H62").Value = 6288.1055
$ws.Range("I62").Value = 5927.4
$ws.Range("J62").Value = 6688.8887
$ws.Range("K62").Value = 5927.4
$ws.Range("L62").Value = 6688.8887
$ws.Range("M62").Value = -5303.4
$ws.Range("N62").Value = -7936.8887
# Row 65
$ws.Range("H65").Value = 6288.1055
$ws.Range("I65").Value = 5927.4
$ws.Range("J65").Value = 6688.8887
$ws.Range("K65").Value = 29637
$ws.Range("L65").Value = 33444.4435
$ws.Range("M65").Value = -26517
$ws.Range("N65").Value = -39684.4435
# Row 70
$ws.Range("H70").Value = 39999
$ws.Range("J70").Value = 39999
$ws.Range("L70").Value = 39999
$ws.Range("N70").Value = -40629
# Row 73
$ws.Range("H73").Value = 39999
$ws.Range("J73").Value = 39999
$ws.Range("L73").Value = 39999
$ws.Range("N73").Value = -42183
# Row 76
$ws.Range("H76").Value = 179212.83
$ws.Range("J76").Value = 179212.83
$ws.Range("L76").Value = 179212.83
$ws.Range("N76").Value = -179842.83
# Row 79
$ws.Range("H79").Value = 179212.83
$ws.Range("J79").Value = 179212.83
$ws.Range("L79").Value = 179212.83
$ws.Range("N79").Value = -181396.83
# Row 81
$ws.Range("H81").Value = 2092.05
$ws.Range("I81").Value = 1867.2354
$ws.Range("K81").Value = 3734.4708
$ws.Range("M81").Value = -2673.4708
# Row 84
$ws.Range("H84").Value = 2092.05
$ws.Range("I84").Value = 1867.2354
$ws.Range("K84").Value = 18672.354
$ws.Range("M84").Value = -13368.354
# Row 132
$ws.Range("H132").Value = 4234.769
$ws.Range("I132").Value = 3105.4
$ws.Range("J132").Value = 7999.3335
$ws.Range("K132").Value = 9316.200000000001
$ws.Range("L132").Value = 23998.0005
$ws.Range("M132").Value = -6786.200000000001
$ws.Range("N132").Value = -29058.0005
